# Logged Week 16 and performed season sim from Week 17
# Update Road ("R") row Target Depth Data for both OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 is the "R" (Road) totals row ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 254   # Short Att
$wsOff.Range("C3").Value = 176   # Short Comp
$wsOff.Range("D3").Value = 59    # Deep Att
$wsOff.Range("E3").Value = 34    # Deep Comp

# --- DEF sheet: row 3 is the "R" (Road) totals row ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 210   # Short Att
$wsDef.Range("C3").Value = 128   # Short Comp
$wsDef.Range("D3").Value = 49    # Deep Att
$wsDef.Range("E3").Value = 18    # Deep Comp
$wsDef.Range("G3").Value = 8     # Deep Int
